# image-showcase.xlsx edit:
#  - add a new "aws.ses" command category to the '#system' sheet, inserted
#    (alphabetically) right after the existing "aws.s3" category.
#  - this means:
#      * a brand-new column is inserted before column C, which becomes the
#        new "aws.ses" column (old column C "base" and everything to its
#        right shifts one column to the right, C:Z -> D:AA);
#      * column A (the "target" list of category names) gets a new entry
#        "aws.ses" inserted at row 3 (right after "aws.s3" at row 2), and
#        every category name that used to occupy rows 3-26 shifts down one
#        row, to rows 4-27;
#      * every defined name in the workbook that pointed into the shifted
#        region is repointed to its new column, and a new defined name
#        "aws.ses" is appended, referring to '#system'!$C$2:$C$3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a new (empty) column before column C. This shifts the
#    existing columns C..Z (every category except "target"/"aws.s3")
#    one column to the right, turning them into D..AA.
# ---------------------------------------------------------------------
$ws.Columns("C").Insert()

# ---------------------------------------------------------------------
# 2) Populate the new column C with the "aws.ses" category: header in
#    row 1, and its two commands in rows 2-3.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "aws.ses"
$ws.Cells.Item(2, 3).Value = "sendMail(profile,to,subject,body)"
$ws.Cells.Item(3, 3).Value = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------
# 3) Column A lists every category name (alphabetically). Insert the new
#    "aws.ses" entry at row 3 (right after "aws.s3" at row 2) by shifting
#    the old rows 3-26 down to rows 4-27, then writing "aws.ses" into the
#    now-vacated row 3. Walk bottom-up so we never overwrite a value
#    before it has been copied down.
# ---------------------------------------------------------------------
$oldTargets = @("base", "csv", "desktop", "excel", "external", "image", "io", "jms", "json", "mail", "number", "pdf", "rdbms", "redis", "sms", "sound", "ssh", "step", "web", "webalert", "webcookie", "ws", "ws.async", "xml")

for ($i = $oldTargets.Length - 1; $i -ge 0; $i--) {
    $destRow = 3 + $i + 1
    $ws.Cells.Item($destRow, 1).Value = $oldTargets[$i]
}
$ws.Cells.Item(3, 1).Value = "aws.ses"

# ---------------------------------------------------------------------
# 4) Repoint every existing defined name whose range moved because of the
#    inserted column, then append the new "aws.ses" name at the end.
#    (date/db/math/mq/nextgen are legacy/unused duplicate names that keep
#    referring to their original columns, unchanged.)
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$D`$2:`$D`$36"
$wb.Names.Item("csv").RefersTo = "='#system'!`$E`$2:`$E`$5"
$wb.Names.Item("desktop").RefersTo = "='#system'!`$F`$2:`$F`$92"
$wb.Names.Item("excel").RefersTo = "='#system'!`$G`$2:`$G`$14"
$wb.Names.Item("external").RefersTo = "='#system'!`$H`$2:`$H`$3"
$wb.Names.Item("image").RefersTo = "='#system'!`$I`$2:`$I`$5"
$wb.Names.Item("io").RefersTo = "='#system'!`$J`$2:`$J`$24"
$wb.Names.Item("jms").RefersTo = "='#system'!`$K`$2:`$K`$4"
$wb.Names.Item("json").RefersTo = "='#system'!`$L`$2:`$L`$14"
$wb.Names.Item("mail").RefersTo = "='#system'!`$M`$2:`$M`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$N`$2:`$N`$15"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$O`$2:`$O`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$P`$2:`$P`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$27"
$wb.Names.Item("web").RefersTo = "='#system'!`$V`$2:`$V`$117"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$Y`$2:`$Y`$17"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AA`$2:`$AA`$11"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$U`$2:`$U`$4"
$wb.Names.Item("redis").RefersTo = "='#system'!`$Q`$2:`$Q`$10"
$wb.Names.Item("sound").RefersTo = "='#system'!`$S`$2:`$S`$5"
$wb.Names.Item("sms").RefersTo = "='#system'!`$R`$2:`$R`$2"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$Z`$2:`$Z`$8"

$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")

# ---------------------------------------------------------------------
# 5) The sheet's recorded dimension was already one column wider than its
#    actual content before this edit (stale "A1:AA117" vs. real content
#    only reaching column Z) and that same one-column-wider bookkeeping
#    carries through the insert (content now reaches AA, dimension should
#    read AB). Touch AB1's number format (no visible change) so the used
#    range keeps extending one column past the real data, like before.
# ---------------------------------------------------------------------
$ws.Range("AB1").NumberFormat = "General"

Write-Host "done"
